$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shared by Overview!E2, Overview!F2, zh-cn!C2, de-de!C2)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn row 2: fill in the "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" columns now that the handback happened.
# ---------------------------------------------------------------------------
$mdFileName = "7c7dbad5-b25f-4e7b-8aa2-e75b26a1ca2b.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aaf2f2f46760e98afc8aef8105473d1932dbef55/e2e/7c7dbad5-b25f-4e7b-8aa2-e75b26a1ca2b.md"

$wsZhCn.Range("I2").Value = $mdFileName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName)
$wsZhCn.Range("J2").Value = "7c7dbad5-b25f-4e7b-8aa2-e75b26a1ca2b.53b48923ef31ca583390d1e9f5b3aec1b1ab15ac.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-10-09 07:27:26"

# ---------------------------------------------------------------------------
# 3. de-de row 2: same, with a later handback datetime.
# ---------------------------------------------------------------------------
$wsDeDe.Range("I2").Value = $mdFileName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName)
$wsDeDe.Range("J2").Value = "7c7dbad5-b25f-4e7b-8aa2-e75b26a1ca2b.53b48923ef31ca583390d1e9f5b3aec1b1ab15ac.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-10-09 07:27:41"

# ---------------------------------------------------------------------------
# 4. Widen the columns that now hold the longer strings above (mirrors the
#    auto-fit Excel performs when a column's content grows).
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666664   # E: zh-cn status
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666664   # F: de-de status

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666664       # C: Status
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664       # I: Latest Target File
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664      # J: Latest Handback File

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666664       # C: Status
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664       # I: Latest Target File
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664      # J: Latest Handback File
